$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("A8")
$r.Value = "Design fehlt"
$f = $r.Font
$f.Size = 10
Write-Host "Done"
